# Update stack trace line numbers in the embedded exception text to reflect
# the move from M2Doc 3.0.0 to 3.1.0 (method/line numbers shifted), and add a
# new "RunBefores.evaluate" frame that appears in the newer JUnit run.

$d = $word.ActiveDocument

# 1) M2DocEvaluator.caseQuery line 559 -> 586
$d.Content.Find.Execute(
    "M2DocEvaluator.caseQuery(M2DocEvaluator.java:559)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "M2DocEvaluator.caseQuery(M2DocEvaluator.java:586)", 2)

# 2) M2DocEvaluator.doSwitch line 1216 -> 1239 (occurs three times; wdReplaceAll
#    replaces every occurrence in the document in one pass)
$d.Content.Find.Execute(
    "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)", 2)

# 3) M2DocEvaluator.caseBlock line 1425 -> 1464
$d.Content.Find.Execute(
    "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1464)", 2)

# 4) M2DocEvaluator.caseDocumentTemplate line 287 -> 296
$d.Content.Find.Execute(
    "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)", 2)

# 5) M2DocEvaluator.generate line 276 -> 281
$d.Content.Find.Execute(
    "M2DocEvaluator.generate(M2DocEvaluator.java:276)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "M2DocEvaluator.generate(M2DocEvaluator.java:281)", 2)

# 6) M2DocUtils.generate line 694 -> 805
$d.Content.Find.Execute(
    "M2DocUtils.generate(M2DocUtils.java:694)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "M2DocUtils.generate(M2DocUtils.java:805)", 2)

# 7) AbstractTemplatesTestSuite.prepareoutputAndGenerate line 480 -> 511
$d.Content.Find.Execute(
    "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)", 2)

# 8) AbstractTemplatesTestSuite.generation line 389 -> 420
$d.Content.Find.Execute(
    "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:420)", 2)

# 9) Insert a new stack frame line for RunBefores.evaluate right before the
#    (second) RunAfters.evaluate line that directly follows a
#    "ParentRunner$2.evaluate(ParentRunner.java:268)" frame.
$d.Content.Find.Execute(
    "ParentRunner`$2.evaluate(ParentRunner.java:268)" + [char]10 + [char]9 + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ParentRunner`$2.evaluate(ParentRunner.java:268)" + [char]10 + [char]9 + "at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)" + [char]10 + [char]9 + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)",
    2)
